$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("savedEachDay")

# Rename cCoefWaterstress (row 16, column A)
$ws.Cells.Item(16, 1).Value = "cCoefWaterstress"

# Fill in column I for rows 33-36 with "NA"
$ws.Cells.Item(33, 9).Value = "NA"
$ws.Cells.Item(34, 9).Value = "NA"
$ws.Cells.Item(35, 9).Value = "NA"
$ws.Cells.Item(36, 9).Value = "NA"

# Row 37: cPhotoDuration
$ws.Cells.Item(37, 1).Value = "cPhotoDuration"
$ws.Cells.Item(37, 2).Value = "computed"
$ws.Cells.Item(37, 3).Value = "numeric"
$ws.Cells.Item(37, 4).Value = "phenology"
$ws.Cells.Item(37, 5).Value = "h"
$ws.Cells.Item(37, 6).Value = "Photoperiod duration"
$ws.Cells.Item(37, 7).Value = "pp"
$ws.Cells.Item(37, 8).Value = "durée journalière ensoleillement"
$ws.Cells.Item(37, 9).Value = "NA"

# Insert 4 new rows after row 37 (rows 38-41)
$ws.Rows.Item(38).Resize(4).Insert()

# Row 38: sThermalUnite
$ws.Cells.Item(38, 1).Value = "sThermalUnite"
$ws.Cells.Item(38, 2).Value = "computed"
$ws.Cells.Item(38, 3).Value = "numeric"
$ws.Cells.Item(38, 4).Value = "phenology"
$ws.Cells.Item(38, 5).Value = "d"
$ws.Cells.Item(38, 6).Value = "Number of Daily temperature unit"
$ws.Cells.Item(38, 7).Value = "DTU"
$ws.Cells.Item(38, 8).Value = "Nombre de jours dans le stade"
$ws.Cells.Item(38, 9).Value = 0

# Row 39: sBiologicalDay
$ws.Cells.Item(39, 1).Value = "sBiologicalDay"
$ws.Cells.Item(39, 2).Value = "computed"
$ws.Cells.Item(39, 3).Value = "numeric"
$ws.Cells.Item(39, 4).Value = "phenology"
$ws.Cells.Item(39, 5).Value = "d"
$ws.Cells.Item(39, 6).Value = "Number of Biological day per calindar day"
$ws.Cells.Item(39, 7).Value = "bd"
$ws.Cells.Item(39, 8).Value = "Nombre de jours"
$ws.Cells.Item(39, 9).Value = 0

# Row 40: pCriticalPhotoPerdiod
$ws.Cells.Item(40, 1).Value = "pCriticalPhotoPerdiod"
$ws.Cells.Item(40, 2).Value = "parameter"
$ws.Cells.Item(40, 3).Value = "numeric"
$ws.Cells.Item(40, 4).Value = "phenology"
$ws.Cells.Item(40, 5).Value = "h"
$ws.Cells.Item(40, 6).Value = "Critical photoperiod"
$ws.Cells.Item(40, 7).Value = "cpp"
$ws.Cells.Item(40, 8).Value = "Seuil photopériode"
$ws.Cells.Item(40, 9).Value = "NA"

# Row 41: pPhotoPeriodSensitivity
$ws.Cells.Item(41, 1).Value = "pPhotoPeriodSensitivity"
$ws.Cells.Item(41, 2).Value = "parameter"
$ws.Cells.Item(41, 3).Value = "numeric"
$ws.Cells.Item(41, 4).Value = "phenology"
$ws.Cells.Item(41, 5).Value = "-"
$ws.Cells.Item(41, 6).Value = "Photoperiod sensitivity coefficient"
$ws.Cells.Item(41, 7).Value = "ppsen"
$ws.Cells.Item(41, 8).Value = "Sensibilité de la plante à la photopériode"
$ws.Cells.Item(41, 9).Value = "NA"

$wb.Save()
